# feat: add 2022-Q1 data
#
# Insert a new "2022-Q1" worksheet between "2021-Q4" and "总计", populate it
# with the quarter's fund-holding detail, and update the "总计" (totals)
# sheet with a new summary row for the quarter.

$wb = $excel.ActiveWorkbook
$sheetQ4 = $wb.Worksheets.Item("2021-Q4")
$sheetTotalOrig = $wb.Worksheets.Item("总计")

# Duplicate the "总计" sheet (rather than Worksheets.Add, which creates a
# blank sheet with none of the workbook's existing sheetPr/margins/styles)
# so the new sheet inherits the same look & feel, then drop it right after
# "2021-Q4" and rename/repurpose it.
$sheetTotalOrig.Copy($null, $sheetQ4)
$newSheet = $wb.Worksheets.Item("总计 (2)")
$newSheet.Name = "2022-Q1"

# ---- Header row (row 1): 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名 ----
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"
# Stamp the (already-bold/centered/bordered) header style across the whole
# B1:H1 span -- E1:H1 are brand new cells that don't have it yet.
$newSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# ---- Data rows: B-G look numeric but are stored as TEXT; only H is a real number ----
# Number format must be switched to Text *before* the values are assigned --
# Excel auto-detects numeric-looking strings at write time, so doing this
# afterwards would not "untype" cells that already parsed as numbers.
$newSheet.Range("B2:G4").NumberFormat = "@"

# Row 2
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "006165"
$newSheet.Cells.Item(2, 3).Value = "建信中证1000指数增强A"
$newSheet.Cells.Item(2, 4).Value = "2.75"
$newSheet.Cells.Item(2, 5).Value = "93.00"
$newSheet.Cells.Item(2, 6).Value = "1.09"
$newSheet.Cells.Item(2, 7).Value = "0.0300"
$newSheet.Cells.Item(2, 8).Value = 4

# Row 3
$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "006166"
$newSheet.Cells.Item(3, 3).Value = "建信中证1000指数增强C"
$newSheet.Cells.Item(3, 4).Value = "0.65"
$newSheet.Cells.Item(3, 5).Value = "93.00"
$newSheet.Cells.Item(3, 6).Value = "1.09"
$newSheet.Cells.Item(3, 7).Value = "0.0071"
$newSheet.Cells.Item(3, 8).Value = 4

# Row 4
$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(4, 2).Value = "013442"
$newSheet.Cells.Item(4, 3).Value = "建信中证1000指数增强E"
$newSheet.Cells.Item(4, 4).Value = "0.02"
$newSheet.Cells.Item(4, 5).Value = "93.00"
$newSheet.Cells.Item(4, 6).Value = "1.09"
$newSheet.Cells.Item(4, 7).Value = "0.0002"
$newSheet.Cells.Item(4, 8).Value = 4

# Column A (row-index cells) carries the same bold/centered/bordered style;
# broadcast it down through the two new rows (3 and 4).
$newSheet.Range("A2").Copy()
$newSheet.Range("A2:A4").PasteSpecial(-4122)

# ---- "总计" sheet: insert a new row above "2021-Q4" for the new quarter ----
$sheetTotal = $wb.Worksheets.Item("总计")
$sheetTotal.Rows.Item(2).Insert()
# Row-insert picks up formatting from the row above (bold header font) --
# strip it back to plain/unstyled, matching the other data rows.
$sheetTotal.Range("B2:D2").ClearFormats()

$sheetTotal.Cells.Item(2, 1).Value = 0
$sheetTotal.Cells.Item(2, 2).Value = "2022-Q1"
$sheetTotal.Cells.Item(2, 3).Value = 3
$sheetTotal.Cells.Item(2, 4).Value = 0.04
$sheetTotal.Range("A3").Copy()
$sheetTotal.Range("A2").PasteSpecial(-4122)

$sheetTotal.Cells.Item(3, 1).Value = 1
